# Update the document per the commit: refresh the date heading and all
# the division problems/answers in the 5-column practice table.
#
# Because several cells share the exact same text (e.g. two cells both
# read "22÷2=11, 0" but must become two different values), we scope each
# Find/Replace to the specific table cell's character range rather than
# doing a single document-wide replace, so each occurrence is updated
# independently and correctly.

$d = $word.ActiveDocument

function Replace-InCell($table, $row, $col, $oldText, $newText) {
    $cell = $table.Cell($row, $col)
    $start = $cell.Range.Start
    $end = $cell.Range.End
    $scoped = $d.Range($start, $end)
    $scoped.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                          $true, 0, $false, $newText, 2) | Out-Null
}

# Date heading above the table.
$d.Content.Find.Execute("2023-10-24 Tuesday", $true, $false, $false, $false, `
                         $false, $true, 1, $false, "2023-10-25 Wednesday", 2) | Out-Null

$tbl = $d.Tables(1)

# Row 1
Replace-InCell $tbl 1 1 "56÷4=14, 0" "34÷4=8, 2"
Replace-InCell $tbl 1 2 "16÷6=2, 4" "68÷5=13, 3"
Replace-InCell $tbl 1 3 "14÷5=2, 4" "67÷4=16, 3"
Replace-InCell $tbl 1 4 "22÷2=11, 0" "24÷6=4, 0"
Replace-InCell $tbl 1 5 "22÷2=11, 0" "21÷2=10, 1"

# Row 5
Replace-InCell $tbl 5 1 "98÷7=14, 0" "92÷7=13, 1"
Replace-InCell $tbl 5 2 "78÷5=15, 3" "28÷9=3, 1"
Replace-InCell $tbl 5 3 "51÷9=5, 6" "54÷7=7, 5"
Replace-InCell $tbl 5 4 "50÷6=8, 2" "55÷8=6, 7"
Replace-InCell $tbl 5 5 "75÷4=18, 3" "57÷9=6, 3"

# Row 9
Replace-InCell $tbl 9 1 "87÷8=10, 7" "39÷2=19, 1"
Replace-InCell $tbl 9 2 "23÷7=3, 2" "89÷4=22, 1"
Replace-InCell $tbl 9 3 "84÷3=28, 0" "17÷8=2, 1"
Replace-InCell $tbl 9 4 "58÷6=9, 4" "29÷4=7, 1"
Replace-InCell $tbl 9 5 "47÷6=7, 5" "54÷3=18, 0"

# Row 13
Replace-InCell $tbl 13 1 "54÷9=6, 0" "27÷8=3, 3"
Replace-InCell $tbl 13 2 "11÷4=2, 3" "48÷7=6, 6"
Replace-InCell $tbl 13 3 "42÷9=4, 6" "63÷2=31, 1"
Replace-InCell $tbl 13 4 "38÷7=5, 3" "73÷7=10, 3"
Replace-InCell $tbl 13 5 "71÷2=35, 1" "28÷3=9, 1"

# Row 17
Replace-InCell $tbl 17 1 "66÷3=22, 0" "68÷9=7, 5"
Replace-InCell $tbl 17 2 "79÷5=15, 4" "28÷7=4, 0"
Replace-InCell $tbl 17 3 "10÷8=1, 2" "85÷8=10, 5"
Replace-InCell $tbl 17 4 "48÷5=9, 3" "35÷3=11, 2"
Replace-InCell $tbl 17 5 "69÷9=7, 6" "74÷9=8, 2"

Write-Output "Done."
